# Auto commit - 08061651
# Adds a new maintenance-report row (row 25) below the existing data, fixes
# the wrap-text styling on row 24's P/AC cells (which had been left out of
# the normal wrap-text formatting used by every other data row), and
# extends the print area / dimension to cover the new row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Seed row 25 with the same cell formatting used throughout the table
#    (copy the previous data row's formats down), then fill in the values.
# ---------------------------------------------------------------------
$ws.Range("A23:AK23").Copy()
$ws.Range("A25:AK25").PasteSpecial(-4122)

$ws.Range("A25").Value = 23
$ws.Range("B25").Value = "維修"
$ws.Range("C25").Value = 2025080650
# Leading apostrophe forces this all-digit customer case number to be
# stored as text instead of being auto-converted to a number.
$ws.Range("D25").Value = "'12014114080501"
$ws.Range("E25").Value = "一般件"
$ws.Range("F25").Value = 2014
$ws.Range("G25").Value = "三重果菜市場"
$ws.Range("H25").Value = "新北市三重區"
$ws.Range("I25").Value = "2025-08-05 14:04:16"
$ws.Range("J25").Value = "星期二"
$ws.Range("K25").Value = "下午"
$ws.Range("L25").Value = "HL25"
$ws.Range("M25").Value = "HL-SC螢幕"
$ws.Range("N25").Value = 2501
$ws.Range("O25").Value = "螢幕畫面閃爍頻繁或無畫面"
$ws.Range("P25").Value = "門市反應SC螢幕(LCD)黑屏顯示無訊號，PING1有通可VNC，門市已嘗試將後方線路重新拔插仍異常，門市告知非與監視器共用螢幕....須請台芝到店協助"
$ws.Range("Q25").Value = "THILF02014"
$ws.Range("R25").Value = "新北一"
$ws.Range("S25").Value = "吳宗鴻"
$ws.Range("T25").Value = 1
$ws.Range("U25").Value = "已完工"
$ws.Range("V25").Value = "2025-08-05 14:06:46"
$ws.Range("W25").Value = "2025-08-06 09:40:00"
$ws.Range("X25").Value = "2025-08-06 11:30:00"
$ws.Range("Y25").Value = "2025-08-06 18:06:00"
$ws.Range("Z25").Value = 1.8
$ws.Range("AA25").Value = ""
$ws.Range("AB25").Value = "到場處理"
$ws.Range("AC25").Value = "螢幕測試正常，SC主機的VGA 孔位損壞，需更換主機，現場的螢幕切換器也損壞無法切換，已告知店員報修，目前將SC和螢幕直接對接使用`n更換SC主機`n換下8114003252`n換上8114004371"
$ws.Range("AD25").Value = ""
$ws.Range("AE25").Value = ""
$ws.Range("AF25").Value = ""
$ws.Range("AG25").Value = ""
$ws.Range("AH25").Value = ""
$ws.Range("AI25").Value = ""
$ws.Range("AJ25").Value = ""
$ws.Range("AK25").Value = "O"

# Re-fit the row height back to the sheet's standard height now that the
# long, wrapped descriptions have been entered (matches how every other
# row in the table is stored -- no explicit/custom row height).
$ws.Rows.Item(25).EntireRow.AutoFit()

# ---------------------------------------------------------------------
# 2. Row 24 was missing the wrap-text formatting applied to every other
#    data row's "report description" / "work content" columns (P, AC).
#    Bring it in line with the rest of the table.
# ---------------------------------------------------------------------
$ws.Range("P24").WrapText = $true
$ws.Range("AC24").WrapText = $true

# ---------------------------------------------------------------------
# 3. Extend the printed area to include the newly added row and move the
#    active selection to the first cell of that row, matching the saved
#    workbook state.
# ---------------------------------------------------------------------
$ws.PageSetup.PrintArea = "A1:AK25"
[void]$ws.Range("A25").Select()

Write-Host "Row 25 added; row 24 wrap formatting fixed; print area extended."
